# Weekly update: add a new "Ajo" (garlic) price record as the new row 228,
# pushing the existing historical records (old rows 228..355) down by one row
# (to new rows 229..356), exactly mirroring how the source system appends the
# latest week's observation at the top of the date-ordered block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 228; this shifts every row from
# 228 downward to 229 downward, and extends the sheet from 355 to 356 rows.
$ws.Rows("228:228").Insert()

# Populate the new row 228 with the latest weekly record.
$ws.Range("A228").Value = 5
$ws.Range("B228").Value = "Macroferia Regional de Talca"
$ws.Range("C228").Value = "Maule"
$ws.Range("D228").Value = 44806
$ws.Range("E228").Value = 7
$ws.Range("F228").Value = 100112003
$ws.Range("G228").Value = "Ajo"
$ws.Range("H228").Value = "Chino"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 300
$ws.Range("K228").Value = 27000
$ws.Range("L228").Value = 27000
$ws.Range("M228").Value = 27000
$ws.Range("N228").Value = "`$/malla 10 kilos"
$ws.Range("O228").Value = "China"
$ws.Range("P228").Value = 2700
$ws.Range("Q228").Value = 10
$ws.Range("R228").Value = "Hortaliza"
